$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.404.87"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "2.646.27"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "2.644.27"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +7.30%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").Value = "3.126.49"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "68.237.76"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "2.645.83"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "363.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.23%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "574.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.00%  "
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.21%  "
$ws.Range("E35").Value = "  +3.13%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("B43").Value = "BabyDogeCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D43").Value = "0.0₆0338"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("E51").Value = "  +1.46%  "
